# Add Table S6: "Number of core OTUs in each group and shared between groups."
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Table S6"

# --- Title (row 1) : bold, same look as the other table titles ---
$srcTitle = $wb.Worksheets.Item("Table S4").Range("A1")
$srcTitle.Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Table S6. Number of core OTUs in each group and shared between groups."

# --- Column headers (row 3): bold + bottom border, same look as other table headers ---
$srcHeader = $wb.Worksheets.Item("Table S4").Range("A3")
$srcHeader.Copy()
$ws.Range("B3:F3").PasteSpecial(-4122)

$headers = @("Summer Wild", "Summer Lab", "Torpor", "IBA", "Spring")
$headerCols = @("B", "C", "D", "E", "F")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($headerCols[$i] + "3").Value = $headers[$i]
}

# --- Row labels (column A, rows 4-8): bold + right border ---
$rowLabels = @("Summer Wild", "Summer Lab", "Torpor", "IBA", "Spring")
$ws.Range("A4:A8").Font.Bold = $true
$ws.Range("A4:A8").Borders.Item(10).LineStyle = 1
for ($i = 0; $i -lt $rowLabels.Length; $i++) {
    $r = 4 + $i
    $ws.Range("A" + $r).Value = $rowLabels[$i]
}

# --- Matrix values (counts of shared core OTUs) ---
$ws.Range("B4").Value = 237
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 47

$ws.Range("C5").Value = 28
$ws.Range("D5").Value = 17
$ws.Range("E5").Value = 19
$ws.Range("F5").Value = 20

$ws.Range("D6").Value = 47
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 31

$ws.Range("E7").Value = 48
$ws.Range("F7").Value = 27

$ws.Range("F8").Value = 112

# --- Leftover formatted (blank) placeholder cells in column A ---
$ws.Range("A12:A18").VerticalAlignment = -4108
$ws.Range("A26:A29").VerticalAlignment = -4108

# --- Column widths (approximate best-fit widths from the source workbook) ---
$ws.Columns.Item(1).ColumnWidth = 12.333333
$ws.Columns.Item(2).ColumnWidth = 12.333333
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 6
$ws.Columns.Item(5).ColumnWidth = 6.666667
$ws.Columns.Item(6).ColumnWidth = 5.666667

# --- Selection / view state matching the source ---
$ws.Range("I15").Select()
$ws.Activate()

Write-Host "Table S6 created"
